# Updated symbol list on Fri Feb 10 23:25:03 UTC 2023 with GitHub Actions
# Refreshes Price / Volume(1h) figures (and the swapped MXToken/BTSEToken rows)
# on the "cryptos" worksheet to match the latest coinranking.com scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'305.82"
$ws.Range("E2").Value = "'-0.23%"

# Row 3
$ws.Range("D3").Value = "'40.67"
$ws.Range("E3").Value = "'2.66%"

# Row 4
$ws.Range("D4").Value = "'5.115"
$ws.Range("E4").Value = "'2.23%"

# Row 5
$ws.Range("D5").Value = "'0.07588"
$ws.Range("E5").Value = "'-1.93%"

# Row 6
$ws.Range("D6").Value = "'4.272"
$ws.Range("E6").Value = "'-0.46%"

# Row 7
$ws.Range("D7").Value = "'1.617"
$ws.Range("E7").Value = "'2.42%"

# Row 8
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.459"
$ws.Range("E8").Value = "'-3.91%"

# Row 9
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9060"
$ws.Range("E9").Value = "'-1.31%"

# Row 10
$ws.Range("D10").Value = "'0.1011"
$ws.Range("E10").Value = "'0.64%"

# Row 11
$ws.Range("D11").Value = "'0.1755"
$ws.Range("E11").Value = "'1.75%"

# Row 12
$ws.Range("D12").Value = "'0.09047"
$ws.Range("E12").Value = "'1.63%"

# Row 13
$ws.Range("D13").Value = "'0.04284"
$ws.Range("E13").Value = "'-2.54%"

# Row 14
$ws.Range("D14").Value = "'0.1055"
$ws.Range("E14").Value = "'-0.34%"

# Row 15
$ws.Range("D15").Value = "'0.001252"
$ws.Range("E15").Value = "'-2.39%"

# Row 16
$ws.Range("D16").Value = "'0.005848"
$ws.Range("E16").Value = "'3.45%"

# Row 17
$ws.Range("D17").Value = "'3.350"
$ws.Range("E17").Value = "'-0.44%"

# Row 18
$ws.Range("E18").Value = "'-2.78%"

# Row 19
$ws.Range("D19").Value = "'6.602"
$ws.Range("E19").Value = "'-6.12%"

# Row 20
$ws.Range("E20").Value = "'-0.59%"

# Row 21
$ws.Range("E21").Value = "'-1.73%"

# Row 22
$ws.Range("E22").Value = "'0.80%"

# Row 23
$ws.Range("D23").Value = "'0.001228"
$ws.Range("E23").Value = "'2.08%"

# Row 24
$ws.Range("D24").Value = "'0.004057"
$ws.Range("E24").Value = "'-0.61%"

# Row 25
$ws.Range("D25").Value = "'0.0001301"
$ws.Range("E25").Value = "'6.13%"

# Row 26
$ws.Range("D26").Value = "'0.0003009"
$ws.Range("E26").Value = "'0.55%"

# Row 38
$ws.Range("E38").Value = "'0.42%"

# Row 39
$ws.Range("D39").Value = "'0.05136"
$ws.Range("E39").Value = "'0.37%"

# Row 40
$ws.Range("D40").Value = "'0.007785"
$ws.Range("E40").Value = "'-2.31%"

# Row 41
$ws.Range("D41").Value = "'0.1295"
$ws.Range("E41").Value = "'-2.42%"

# Row 42
$ws.Range("D42").Value = "'0.007064"
$ws.Range("E42").Value = "'-4.42%"

# Row 43
$ws.Range("D43").Value = "'0.001919"
$ws.Range("E43").Value = "'-3.83%"

# Row 44
$ws.Range("D44").Value = "'0.008440"
$ws.Range("E44").Value = "'4.87%"

# Row 45
$ws.Range("D45").Value = "'0.3330"
$ws.Range("E45").Value = "'0.10%"

# Row 46
$ws.Range("D46").Value = "'0.00006364"
$ws.Range("E46").Value = "'-4.76%"

# Row 47
$ws.Range("E47").Value = "'-0.39%"

# Row 48
$ws.Range("D48").Value = "'0.004403"
$ws.Range("E48").Value = "'6.89%"

# Row 49
$ws.Range("D49").Value = "'0.006999"
$ws.Range("E49").Value = "'110.28%"

# Row 50
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'-0.39%"

# Row 51
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'-0.39%"
